$wb = $excel.ActiveWorkbook

# --- Sheet "Đơn sale chính": insert a new detail row (row 3) and push
# the existing "Tổng" (total) row down to row 4, updating its totals.
$ws1 = $wb.Worksheets.Item("Đơn sale chính")

$ws1.Rows.Item(3).Insert()

$ws1.Range("A3").Value = "HD-LUXURY"
$ws1.Range("B3").Value = 643
# Keep the date as literal text (matches existing column formatting)
# instead of letting it be auto-parsed into a date serial number.
$ws1.Range("C3").NumberFormat = "@"
$ws1.Range("C3").Value = "08-11-2024"
$ws1.Range("C3").ClearFormats()
$ws1.Range("D3").Value = "CẦN THƠ"
$ws1.Range("E3").Value = "Trần Thị Nhi"
$ws1.Range("F3").Value = "Cá nhân"
$ws1.Range("G3").Value = "Cắt mí"
$ws1.Range("H3").Value = 3000000
$ws1.Range("I3").Value = 0
$ws1.Range("J3").Value = 0
$ws1.Range("K3").Value = 3000000
$ws1.Range("L3").Value = 3000000
$ws1.Range("M3").Value = 0.1
$ws1.Range("N3").Value = 300000

# Update the "Tổng" row (now row 4) with the new aggregated totals.
$ws1.Range("B4").Value = 2
$ws1.Range("H4").Value = 5100000
$ws1.Range("K4").Value = 5100000
$ws1.Range("L4").Value = 5100000
$ws1.Range("N4").Value = 510000

# --- Sheet "Lương": refresh the computed payroll summary figures.
$ws4 = $wb.Worksheets.Item("Lương")

$ws4.Range("B2").Value = 9.5
$ws4.Range("B3").Value = 332500
$ws4.Range("B4").Value = 1696428.571428572
$ws4.Range("B5").Value = 510000
$ws4.Range("B35").Value = 2848928.571428572
$ws4.Range("B38").Value = 2848928.571428572
